$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 356.0625
$ws.Range("I2").Value = 95.25
$ws.Range("J2").Value = 443
$ws.Range("K2").Value = 95.25
$ws.Range("L2").Value = 443
$ws.Range("M2").Value = 17.75
$ws.Range("N2").Value = -669

$ws.Range("H13").Value = 305
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = $null

$ws.Range("H40").Value = 4210.5835
$ws.Range("I40").Value = 5465.6
$ws.Range("J40").Value = 3314.1428
$ws.Range("K40").Value = 5465.6
$ws.Range("L40").Value = 3314.1428
$ws.Range("M40").Value = -5290.6
$ws.Range("N40").Value = -3664.1428

$ws.Range("H41").Value = 1374.75
$ws.Range("I41").Value = 1957.5714
$ws.Range("K41").Value = 1957.5714
$ws.Range("M41").Value = -1517.5714

$ws.Range("H43").Value = 1360.9166
$ws.Range("I43").Value = 1080.6666
$ws.Range("J43").Value = 1454.3334
$ws.Range("K43").Value = 1080.6666
$ws.Range("L43").Value = 1454.3334
$ws.Range("M43").Value = -1011.6666
$ws.Range("N43").Value = -1592.3334

$ws.Range("H51").Value = 13892248
$ws.Range("I51").Value = 7000
$ws.Range("J51").Value = 15875854
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 15875854
$ws.Range("M51").Value = -6516
$ws.Range("N51").Value = -15876822

$ws.Range("H76").Value = 79644.92999999999
$ws.Range("J76").Value = 3490.9167
$ws.Range("L76").Value = 3490.9167
$ws.Range("N76").Value = -4120.9167

$ws.Range("H79").Value = 79644.92999999999
$ws.Range("J79").Value = 3490.9167
$ws.Range("L79").Value = 3490.9167
$ws.Range("N79").Value = -5674.9167

$ws.Range("H93").Value = 37601
$ws.Range("J93").Value = 37601
$ws.Range("L93").Value = 37601
$ws.Range("N93").Value = -42593

$ws.Range("H113").Value = 2380.7896
$ws.Range("I113").Value = 1773.5
$ws.Range("J113").Value = 3055.5557
$ws.Range("K113").Value = 1773.5
$ws.Range("L113").Value = 3055.5557
$ws.Range("M113").Value = 1480.5
$ws.Range("N113").Value = -9563.555700000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3814.1667
$ws.Range("J63").Value = 4099
$ws.Range("L63").Value = 4099
$ws.Range("N63").Value = -5471

$ws.Range("H66").Value = 3814.1667
$ws.Range("J66").Value = 4099
$ws.Range("L66").Value = 20495
$ws.Range("N66").Value = -27359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null

$ws.Range("H63").Value = 15500
$ws.Range("J63").Value = 15500
$ws.Range("L63").Value = 15500
$ws.Range("N63").Value = -16872

$ws.Range("H66").Value = 15500
$ws.Range("J66").Value = 15500
$ws.Range("L66").Value = 46500
$ws.Range("N66").Value = -53364

$ws.Range("H99").Value = 3973933
$ws.Range("I99").Value = 5958766
$ws.Range("J99").Value = 4266.6665
$ws.Range("K99").Value = 5958766
$ws.Range("L99").Value = 4266.6665
$ws.Range("M99").Value = -5957268
$ws.Range("N99").Value = -7262.6665

$ws.Range("H126").Value = 3973933
$ws.Range("I126").Value = 5958766
$ws.Range("J126").Value = 4266.6665
$ws.Range("K126").Value = 17876298
$ws.Range("L126").Value = 12799.9995
$ws.Range("M126").Value = -17873828
$ws.Range("N126").Value = -17739.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1453
$ws.Range("I17").Value = 131.5
$ws.Range("J17").Value = 2334
$ws.Range("K17").Value = 394.5
$ws.Range("L17").Value = 7002
$ws.Range("M17").Value = -225.5
$ws.Range("N17").Value = -7340

$ws.Range("H42").Value = 8166.6665
$ws.Range("J42").Value = 8166.6665
$ws.Range("L42").Value = 24499.9995
$ws.Range("N42").Value = -25567.9995

$ws.Range("H80").Value = 966.38464
$ws.Range("I80").Value = 850
$ws.Range("J80").Value = 976.0833
$ws.Range("K80").Value = 2550
$ws.Range("L80").Value = 2928.2499
$ws.Range("M80").Value = -1614
$ws.Range("N80").Value = -4800.2499

$ws.Range("H83").Value = 966.38464
$ws.Range("I83").Value = 850
$ws.Range("J83").Value = 976.0833
$ws.Range("K83").Value = 7650
$ws.Range("L83").Value = 8784.7497
$ws.Range("M83").Value = -2970
$ws.Range("N83").Value = -18144.7497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3919
$ws.Range("I80").Value = 7276.25
$ws.Range("J80").Value = 2698.182
$ws.Range("K80").Value = 7276.25
$ws.Range("L80").Value = 2698.182
$ws.Range("M80").Value = -6278.25
$ws.Range("N80").Value = -4694.182

$ws.Range("H83").Value = 3919
$ws.Range("I83").Value = 7276.25
$ws.Range("J83").Value = 2698.182
$ws.Range("K83").Value = 36381.25
$ws.Range("L83").Value = 13490.91
$ws.Range("M83").Value = -31389.25
$ws.Range("N83").Value = -23474.91

$ws.Range("H102").Value = 1219.5714
$ws.Range("I102").Value = 1055.4
$ws.Range("J102").Value = 1630
$ws.Range("K102").Value = 1055.4
$ws.Range("L102").Value = 1630
$ws.Range("M102").Value = 566.5999999999999
$ws.Range("N102").Value = -4874

$ws.Range("H126").Value = 3409.4119
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 4662.222
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 13986.666
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -18926.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1620.9706
$ws.Range("I7").Value = 1265.1765
$ws.Range("J7").Value = 1976.7646
$ws.Range("K7").Value = 1265.1765
$ws.Range("L7").Value = 1976.7646
$ws.Range("M7").Value = -1153.1765
$ws.Range("N7").Value = -2200.7646

$ws.Range("H16").Value = 966.6667
$ws.Range("J16").Value = 1025
$ws.Range("L16").Value = 1025
$ws.Range("N16").Value = -1365

$ws.Range("H46").Value = 3074.75
$ws.Range("I46").Value = 5349.5
$ws.Range("J46").Value = 800
$ws.Range("K46").Value = 5349.5
$ws.Range("L46").Value = 800
$ws.Range("M46").Value = -5161.5
$ws.Range("N46").Value = -1176

$ws.Range("H55").Value = 212.33333
$ws.Range("I55").Value = 233.7619
$ws.Range("K55").Value = 233.7619
$ws.Range("M55").Value = -60.7619

$ws.Range("H68").Value = 33836400
$ws.Range("I68").Value = 112778296
$ws.Range("J68").Value = 4158.5713
$ws.Range("K68").Value = 112778296
$ws.Range("L68").Value = 4158.5713
$ws.Range("M68").Value = -112777547
$ws.Range("N68").Value = -5656.5713

$ws.Range("H71").Value = 33836400
$ws.Range("I71").Value = 112778296
$ws.Range("J71").Value = 4158.5713
$ws.Range("K71").Value = 563891480
$ws.Range("L71").Value = 20792.8565
$ws.Range("M71").Value = -563887736
$ws.Range("N71").Value = -28280.8565

$ws.Range("H126").Value = 1620.9706
$ws.Range("I126").Value = 1265.1765
$ws.Range("J126").Value = 1976.7646
$ws.Range("K126").Value = 3795.5295
$ws.Range("L126").Value = 5930.293799999999
$ws.Range("M126").Value = -1325.5295
$ws.Range("N126").Value = -10870.2938

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5089.28
$ws.Range("I62").Value = 4839.619
$ws.Range("J62").Value = 6400
$ws.Range("K62").Value = 4839.619
$ws.Range("L62").Value = 6400
$ws.Range("M62").Value = -4215.619
$ws.Range("N62").Value = -7648

$ws.Range("H65").Value = 5089.28
$ws.Range("I65").Value = 4839.619
$ws.Range("J65").Value = 6400
$ws.Range("K65").Value = 24198.095
$ws.Range("L65").Value = 32000
$ws.Range("M65").Value = -21078.095
$ws.Range("N65").Value = -38240

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null

$ws.Range("H113").Value = 447.6842
$ws.Range("I113").Value = 325.2143
$ws.Range("J113").Value = 790.6
$ws.Range("K113").Value = 975.6428999999999
$ws.Range("L113").Value = 2371.8
$ws.Range("M113").Value = 1194.3571
$ws.Range("N113").Value = -6711.8

$ws.Range("H126").Value = 2609.45
$ws.Range("I126").Value = 3125.9333
$ws.Range("J126").Value = 1060
$ws.Range("K126").Value = 9377.7999
$ws.Range("L126").Value = 3180
$ws.Range("M126").Value = -6907.7999
$ws.Range("N126").Value = -8120
